$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (I2): fix wording 'Enter Y or N .' -> 'Enter Y or N.'
$textI2 = @'
While we could continue with that, let's move to a different question. Do Framingham risk estimates correlate with an individual's annual medical expenditures? It would make sense if they did- a high risk person might be engaging more with their doctors to begin with, or a person with lots of medical expenses might have more stress and less means/bandwidth to sustain a healthy diet, inflating some parts of the Framingham scores (especially systolic blood Enterure and cholesterol).

Should I run that analysis? Enter Y or N.
'@
$ws.Range("I2").Value = $textI2

# Column I width for the new wider text column
$ws.Columns.Item(9).ColumnWidth = 41.4

# New row 9
$textI9 = @'
Enter 1 for the same model with standardized coefficients 
Enter 2 for a review of the regression diagnostics
Enter 3 for a full interpretation of the model and diagnostics
Enter 4 to move on.
All options except for #4 will return you to these options so you can see everything. 
'@
$ws.Range("I9").Value = $textI9
$ws.Range("I9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 409.5

# New row 10
$textI10 = @'
MODEL FIT:
The model's adjusted R-squared value is .389, meaning that these three variables can collectively explain 38.7% of the variance in annual medical spending. Though certainly not guaranteed accuracy, this is a really high number for social sciences data and a model with only two variables!
COEFFICIENT ESTIMATES:
The intercept reflects all covariates=0, reflecting the estimated annual spending for a female with 0 estimated risk. In this case, that's still $2,327, still rather high!
Risk_pct's estimate is statistically significant (dubious use for inference in this case), indicating that the result probably shouldn't be chocked up to random luck of the draw. The coefficient tells us that as a person's risk of cardiovascular disease increases by 1%, we expect that their annual medical spending will increase by $54. Though that sounds small, recall that the risk scores go up to 30, so that's as much as a 54*30=$1,620 increase. Be mindful, though, that this is a description of our current dataset's patterns, not necessarily a prediction of the future.
The coefficient for gender is also statistically significant, estimating that  males in these data spend about $748 less than females on medical costs per year. 
The coefficient for SNAP is small and is not statistically significant- it should proably be ignored. Notice the size of the standard error, which is almost equal to the size of the coefficient estimate- we can't conclude if SNAP membership is associated with an increase of decrease in data.  
If we were to look at the standardized coefficients, we can confirm that the correlation for risk scores is the strongest in the model, with a 1 std deviation increase in risk % corresponding with a .59 standard deviation increase in medical spending. Contrast that with the -.487 standard deviation change when comparing males to females.
DIAGNOSTICS:
The residual plots give some pause- notice how in some spots the residuals are tightly clustered around the reference line, which in other spots they are more spread out. This is an indication of heteroskedasticity, though the Breusch Pagan test is not significant, suggesting heteroskedasticity is not a massive problem. But the CCPR plot again suggests a shape that is not entirely linear, with a bit of a curve visible. This, and the histogram of errors that is clearly not normal, suggests that this model is still not a great fit to the data. We could try to fix this by introducing more covariates, if possible, or trying different functional forms for the variables we've already got, such as a quadratic term for the Framingham risk pct.
The outlier plot indicates six cases that are likely to be outliers. We should run a later model without those cases to see if they are influencing our conclusions at all. 
'@
$ws.Range("I10").Value = $textI10
$ws.Range("I10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 409.5

# New row 11
$textI11 = @'
The square term considerably improves the models fit a bit, as evinced by the increased r-squared. But it's also much harder to interpret what is going on now because the relationship between the risk scores and annual medical spending is now curvilinear. We need to visualize the estimates to help us see what is going on.
Press enter to graph the results.
'@
$ws.Range("I11").Value = $textI11
$ws.Range("I11").WrapText = $true
$ws.Rows.Item(11).RowHeight = 130.5

# New row 12
$textI12 = @'
Looking at this graph, we can see that spending is expected to increase sharply for about the first 14% in the risk scale, at which point annual spending tapers out. We can use this regression equation to make particular estimates of a person's expenditures based on their characteristics if we wanted to- just remember to factor gender in, too!
That concludes the regression analysis- this is probably the best fitting model we can get with such limited data and a basic linear regression. However, we should be cautious about trusting this model- the omission of many key variables is likely to have biased out models in ways that we cannot know. 
******Press enter to go back to the main hub.******
'@
$ws.Range("I12").Value = $textI12
$ws.Range("I12").WrapText = $true
$ws.Rows.Item(12).RowHeight = 275.5

# Update selection to match the new active cell
$ws.Range("I11").Select() | Out-Null
